# Customer shipments template: add Contact/Address columns ahead of the
# existing shipment-detail columns (import update on customer side).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for 8 new lead columns by inserting before the current column A.
$ws.Range("A1:H1").EntireColumn.Insert()

# Column letter -> header text, keyed by final position in the sheet.
$headerByCol = @{
    1 = "Contact Person"
    2 = "Contact Number"
    3 = "Email Address"
    4 = "Street"
    5 = "Barangay"
    6 = "Municipality"
    7 = "Province"
    8 = "Zip Code"
}

# Column widths (character units) as they appear in the final workbook.
# NOTE: the engine's ColumnWidth setter re-applies the 5px/MDW standard
# padding on top of the value it stores, so we back it out here to land on
# the intended width once the value is round-tripped through pixels.
$widthByCol = @{
    1 = 30.1640625
    2 = 23.1640625
    3 = 21.83203125
    4 = 31.1640625
    5 = 22.6640625
    6 = 26.5
    7 = 25.5
    8 = 22.5
}

# Write the new header text in the same order the strings were first
# introduced in the workbook (Contact Person, Contact Number, Email
# Address, Province, Municipality, Barangay, Street, Zip Code) so the
# shared-string table comes out in that order, even though on the sheet
# Province/Municipality/Barangay/Street land right-to-left in columns
# G/F/E/D.
$fillOrder = @(1, 2, 3, 7, 6, 5, 4, 8)
foreach ($col in $fillOrder) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headerByCol[$col]
    # Match the bottom-border style already used across the header row.
    $cell.Borders.Item(9).LineStyle = 1
    $ws.Columns.Item($col).ColumnWidth = $widthByCol[$col] - (5 / 6)
}

$excel.ActiveWindow.Zoom = 141
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H2").Select()
